$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 8, shifting existing rows 8:51 down to 9:52.
$ws.Rows.Item(8).Insert()

# Populate the newly inserted row 8 with the new record.
$ws.Range("A8").Value = 1
$ws.Range("B8").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C8").Value = "Arica y Parinacota"
$ws.Range("D8").Value = 44602
$ws.Range("D8").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E8").Value = 15
$ws.Range("F8").Value = 100112009
$ws.Range("G8").Value = "Acelga"
$ws.Range("H8").Value = "Sin especificar"
$ws.Range("I8").Value = "Primera"
$ws.Range("J8").Value = 200
$ws.Range("K8").Value = 2000
$ws.Range("L8").Value = 2500
$ws.Range("M8").Value = 2250
$ws.Range("N8").Value = "$/atado 2,5 a 3 kilos"
$ws.Range("O8").Value = "Región de Arica y Parinacota"
$ws.Range("P8").Value = 750
$ws.Range("Q8").Value = 3
$ws.Range("R8").Value = "Hortaliza"
